$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (mirror H1's text/style pattern)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style (font/border/alignment) from H1 into I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for columns I and J (rows 2-27)
$values = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(8, 8)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(6, 6)
    8  = @(4, 5)
    9  = @(6, 6)
    10 = @(6, 6)
    11 = @(6, 6)
    12 = @(6, 6)
    13 = @(6, 7)
    14 = @(6, 7)
    15 = @(7, 8)
    16 = @(10, 11)
    17 = @(8, 8)
    18 = @(6, 6)
    19 = @(7, 7)
    20 = @(8, 9)
    21 = @(6, 7)
    22 = @(7, 7)
    23 = @(9, 9)
    24 = @(7, 8)
    25 = @(4, 6)
    26 = @(9, 9)
    27 = @(9, 9)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
